# Applies the "Saldo.xlsx" update:
#  - Remove 5 accounts that dropped out of the export (YURI, LEVI, CASSIO,
#    CAROLINE, PAULO)
#  - Add one new account (PEDRO / 005324840 / 146.5) that now appears in the
#    export, inserted right after account 005022526 (ALEXANDRE / 147.18) so
#    the "Saldo" column stays sorted in descending order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = $ws.Columns(1)

# Locate each row to delete by its account number (column A) and remove it.
# Delete highest row first so the row numbers found for the others stay valid.
$accountsToDelete = @("004612043", "005206566", "004508526", "004221638", "004572740")

$rowsToDelete = @()
foreach ($acct in $accountsToDelete) {
    $cell = $col.Find($acct)
    $rowsToDelete += $cell.Row
}

$rowsToDelete = $rowsToDelete | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $ws.Rows($r).Delete()
}

# Find the row to insert after (account 005022526 / ALEXANDRE / 147.18),
# re-resolved after the deletions above shifted rows upward.
$anchor = $col.Find("005022526")
$insertRow = $anchor.Row + 1

$ws.Rows($insertRow).Insert()

# Force column A to text so the leading zeros in the account number survive.
$ws.Cells.Item($insertRow, 1).NumberFormat = "@"
$ws.Cells.Item($insertRow, 1).Value = "005324840"
$ws.Cells.Item($insertRow, 2).Value = "PEDRO"
$ws.Cells.Item($insertRow, 3).Value = 146.5
